# Update "想去人数" (want-to-go count) figures across the workbook's four
# sheets to match the refreshed scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1236
$ws.Range("F4").Value = 51
$ws.Range("F5").Value = 3468
$ws.Range("F6").Value = 1747
$ws.Range("F7").Value = 6284
$ws.Range("F8").Value = 134
$ws.Range("F9").Value = 1885
$ws.Range("F10").Value = 500
$ws.Range("F12").Value = 26
$ws.Range("F15").Value = 46
$ws.Range("F16").Value = 7447
$ws.Range("F18").Value = 58
$ws.Range("F20").Value = 106
$ws.Range("F21").Value = 1730
$ws.Range("F28").Value = 1670
$ws.Range("F29").Value = 787
$ws.Range("F30").Value = 351
$ws.Range("F34").Value = 86

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 204

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9536
$ws.Range("F5").Value = 258

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9536
$ws.Range("F5").Value = 1236
$ws.Range("F7").Value = 51
$ws.Range("F10").Value = 3469
$ws.Range("F11").Value = 258
$ws.Range("F12").Value = 1747
$ws.Range("F13").Value = 6284
$ws.Range("F14").Value = 1885
$ws.Range("F16").Value = 500
$ws.Range("F18").Value = 26
$ws.Range("F21").Value = 46
$ws.Range("F22").Value = 7447
$ws.Range("F24").Value = 58
$ws.Range("F26").Value = 106
$ws.Range("F27").Value = 1730
$ws.Range("F34").Value = 1670
$ws.Range("F35").Value = 787
$ws.Range("F37").Value = 351
$ws.Range("F42").Value = 86
